# Applies the "Added more architecture diagrams" edit:
#  1. Re-stacks three labelled rectangles (STASHL, MBR, IR) so they are
#     drawn after (on top of) the "RAM/ROM" rectangle, i.e. moved to the
#     end of the slide's z-order, in the order STASHL, MBR, IR.
#  2. Nudges three "Up-Down Arrow" connector shapes to new positions.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# --- 1. Bring the three rectangles to the front, in order -----------------
# (msoBringToFront = 0) moves a shape to the end of the z-order (last
# drawn / topmost), matching where the diff re-inserts these <p:sp>
# elements -- right after the "RAM/ROM" shape at the end of the spTree.

$stashl = Get-ShapeById $s 234   # "Rectangle 233" / text "STASHL"
$mbr    = Get-ShapeById $s 9     # "Rectangle 8"   / text "MBR"
$ir     = Get-ShapeById $s 6     # "Rectangle 5"   / text "IR"

$stashl.ZOrder(0)
$mbr.ZOrder(0)
$ir.ZOrder(0)

# --- 2. Reposition the three up-down arrows --------------------------------

$arrow9   = Get-ShapeById $s 10   # "Up-Down Arrow 9"
$arrow135 = Get-ShapeById $s 136  # "Up-Down Arrow 135"
$arrow142 = Get-ShapeById $s 143  # "Up-Down Arrow 142"

$arrow9.Left   = 319.6353607307101
$arrow9.Top    = 163.64732283464568

$arrow135.Left = 486.5607874015748
$arrow135.Top  = 162.561027622047

$arrow142.Left = 320.59102362204726
$arrow142.Top  = 28.073700987401576
